$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Range("A1").Value = "Datos actualizados a 25 de Abril de 2020 a las 01:22"

# Row 4
$ws.Range("B4").Value = 923474
$ws.Range("C4").Value = 37032
$ws.Range("D4").Value = 93369
$ws.Range("E4").Value = 778020
$ws.Range("F4").Value = 14954
$ws.Range("G4").Value = 1851
$ws.Range("H4").Value = 52085

# Row 85
$ws.Range("A85").Value = "Nigeria"
$ws.Range("B85").Value = 1095
$ws.Range("C85").Value = 114
$ws.Range("D85").Value = 208
$ws.Range("E85").Value = 855
$ws.Range("F85").Value = 2
$ws.Range("G85").Value = 1
$ws.Range("H85").Value = 32

# Row 86
$ws.Range("A86").Value = "Costa de Marfil"
$ws.Range("B86").Value = 1077
$ws.Range("C86").Value = 73
$ws.Range("D86").Value = 419
$ws.Range("E86").Value = 644
$ws.Range("F86").Value = 0
$ws.Range("H86").Value = 14

# Row 87
$ws.Range("A87").Value = "Hong Kong"
$ws.Range("B87").Value = 1036
$ws.Range("C87").Value = 0
$ws.Range("D87").Value = 699
$ws.Range("E87").Value = 333
$ws.Range("F87").Value = 9
$ws.Range("H87").Value = 4

# Row 88
$ws.Range("A88").Value = "Republica de Yibuti"
$ws.Range("B88").Value = 999
$ws.Range("C88").Value = 13
$ws.Range("D88").Value = 330
$ws.Range("E88").Value = 667
$ws.Range("F88").Value = 0
$ws.Range("H88").Value = 2

# Row 102
$ws.Range("A102").Value = "Uruguay"
$ws.Range("B102").Value = 563
$ws.Range("C102").Value = 14
$ws.Range("D102").Value = 369
$ws.Range("E102").Value = 182
$ws.Range("F102").Value = 9
$ws.Range("H102").Value = 12

# Row 103
$ws.Range("A103").Value = "Honduras"
$ws.Range("B103").Value = 562
$ws.Range("C103").Value = 43
$ws.Range("D103").Value = 50
$ws.Range("E103").Value = 465
$ws.Range("H103").Value = 47

# Row 138
$ws.Range("D138").Value = 131
$ws.Range("E138").Value = 2

# Row 145
$ws.Range("B145").Value = 111
$ws.Range("C145").Value = 4
$ws.Range("D145").Value = 87
$ws.Range("E145").Value = 23
$ws.Range("F145").Value = 0

# Row 165
$ws.Range("B165").Value = 52
$ws.Range("C165").Value = 2
$ws.Range("E165").Value = 49

# Row 188
$ws.Range("A188").Value = "Republica de Africa Central"
$ws.Range("D188").Value = 10
$ws.Range("E188").Value = 6

# Row 189
$ws.Range("A189").Value = "Dominica"

# Row 202
$ws.Range("A202").Value = "Surinam"
$ws.Range("D202").Value = 6
$ws.Range("E202").Value = 3

# Row 203
$ws.Range("A203").Value = "Gambia"
$ws.Range("D203").Value = 8
$ws.Range("E203").Value = 1

Write-Host "Edit complete"